$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r, A(date-serial), B(nuovi pos.), C(somma mobile 7gg.), D(somma mobile 7gg. per 100mila ab.)
$data = @(
    @(270, 44344, 0, 1, 6.628222973420826),
    @(271, 44345, 0, 1, 6.628222973420826),
    @(272, 44346, 0, 0, 0),
    @(273, 44347, 0, 0, 0),
    @(274, 44348, 0, 0, 0),
    @(275, 44349, 0, 0, 0),
    @(276, 44350, 1, 1, 6.628222973420826),
    @(277, 44351, 0, 1, 6.628222973420826),
    @(278, 44352, 0, 1, 6.628222973420826),
    @(279, 44353, 0, 1, 6.628222973420826),
    @(280, 44354, 0, 1, 6.628222973420826),
    @(281, 44355, 0, 1, 6.628222973420826),
    @(282, 44356, 0, 1, 6.628222973420826),
    @(283, 44357, 0, 0, 0),
    @(284, 44358, 0, 0, 0),
    @(285, 44359, 1, 1, 6.628222973420826),
    @(286, 44360, 1, 2, 13.25644594684165),
    @(287, 44361, 0, 2, 13.25644594684165),
    @(288, 44362, 1, 3, 19.88466892026248),
    @(289, 44363, 0, 3, 19.88466892026248),
    @(290, 44364, 0, 3, 19.88466892026248),
    @(291, 44365, 0, 3, 19.88466892026248),
    @(292, 44366, 0, 2, 13.25644594684165),
    @(293, 44367, 0, 1, 6.628222973420826),
    @(294, 44368, 0, 1, 6.628222973420826),
    @(295, 44369, 0, 0, 0),
    @(296, 44370, 0, 0, 0),
    @(297, 44371, 0, 0, 0),
    @(298, 44372, 0, 0, 0),
    @(299, 44373, 0, 0, 0),
    @(300, 44374, 0, 0, 0),
    @(301, 44375, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Copy the formatting (style) of the date column from the last existing row
# down onto all of the newly added date cells, matching the existing pattern.
$ws.Range("A269").Copy()
$ws.Range("A270:A301").PasteSpecial(-4122)
